# Insert a new "Supplementary File 1 (...)" list item right after the
# "sample_details_et_mix5_name_id_atom.csv (input file for TSR key
# generation code for pigments using all atoms except hydrogen atoms)"
# bullet, matching the commit's "Add files via upload" change.

$d = $word.ActiveDocument

# Locate the paragraph whose text starts with the
# sample_details_et_mix5_name_id_atom.csv bullet so the insertion point is
# robust even if paragraph indices shift.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("sample_details_et_mix5_name_id_atom.csv")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a brand-new paragraph right after the target bullet; it
    # inherits the target's paragraph formatting (ListParagraph style +
    # the same numbered/bulleted list) the same way Word does when you
    # press Enter at the end of a list item.
    $target.Range.InsertParagraphAfter()

    $newPara = $target.Next()
    $newRange = $newPara.Range
    $newRange.Collapse(1)

    $newRange.InsertAfter("Supplementary File")
    $newRange.Collapse(0)
    $newRange.InsertAfter(" ")
    $newRange.Collapse(0)
    $newRange.InsertAfter("1 (")
    $newRange.Collapse(0)
    $newRange.InsertAfter("Supplementary_File1_Sample_Details.csv")
    $newRange.Collapse(0)
    $newRange.InsertAfter(" for protein global analyses")
    $newRange.Collapse(0)
    $newRange.InsertAfter(")")
}
